$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Rows in the "Ready for handoff" block (031cee3f, 0b142384, 527bef03,
# 7054c5c4, 74956ebe, 87a0e161) whose Priority column (E) needs to be
# set from blank to "ht" on both the zh-cn and de-de sheets.
$priorityRows = @(7, 8, 10, 11, 12, 13)
foreach ($r in $priorityRows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}

# "Latest HO Xliff Generate Date" (Overview!G7) and the shared
# "Latest Handoff Datetime" text on de-de!H7 both displayed
# "2016-08-25 04:20:57" -> bump to "2016-08-25 04:21:18" for every row
# that shared that timestamp string.
$overviewDateRows = @(7, 8, 10, 11, 12, 13)
foreach ($r in $overviewDateRows) {
    $wsOverview.Range("G$r").Value = "2016-08-25 04:21:18"
    $wsDeDe.Range("H$r").Value = "2016-08-25 04:21:18"
}

# zh-cn!H7's own "Latest Handoff Datetime" text
# "2016-08-25 04:20:52" -> "2016-08-25 04:21:13".
$zhCnDateRows = @(7, 8, 10, 11, 12, 13)
foreach ($r in $zhCnDateRows) {
    $wsZhCn.Range("H$r").Value = "2016-08-25 04:21:13"
}
